# PowerShell-style Excel COM-interop script
# Updates the cryptos price/volume table to reflect the latest scrape,
# and reorders two groups of rows (44-46 and 48-49) to match the
# new ranking order, per commit "Updated cryptos list ... GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain plain text (avoid Excel auto-number conversion),
# matching the original inlineStr cell values exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.578.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.382.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.51%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.68'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.36%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.18'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0923'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.56'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.37%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.985'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.744.79'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.56'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.376.34'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.484.65'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.29'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +13.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.79%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.69'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.49'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.25'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.45%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.55%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.51'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.94%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.55'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0959'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.47'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.30'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.88'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.83%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.17%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.27%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.77'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.99'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.03'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0357'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.36'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.11'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.23%  '

# Rows 44-46 and 48-49: three coins were reordered in the new ranking
# (Celestia/Algorand moved above Maker; FirstDigitalUSD moved above ordi),
# with refreshed Price/Volume(1h) values.
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.12'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.27%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.231'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.15%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.879.52'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +13.18%  '

$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.83'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.10'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.90%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.41'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.04'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.99%  '
